$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("od1")

# "Calcium Carbide" -> "calcium carbide" (de-capitalize) in the 2012 JV paragraph
$ws.Range("D15").Value = '<div class="timeline-date">2012</div>ALMAMET formed a Joint Venture with an Indian company in 2001 known as ALMAMET Reagents (India) Pvt Ltd. Through this JV ALMAMET was catering to steel plants e.g. TATA, JSW, BHUSHAN, ESSAR & SAIL. ALMAMET also had a production unit of calcium carbide blends in Bhilai. In 2012 ALMAMET decides to run its own company by establishing ALMAMET India Private Limited. Since the foundation in 2012 ALMAMET India Private Limited is in discussion with the Indian steel industry to invest in various projects. At the same time ALMAMET India Private Limited is in discussion with companies from chemical & pyrotechnical sector. In near future the ALMAMET group will take major investments to increase the representation within the subcontinent of India. '

# "Warehouse" -> "warehouse"
$ws.Range("C18").Value = '<div class="timeline-title">Almamet India warehouse, Karnataka.</div>'

# "Symposium" -> "symposium"
$ws.Range("C20").Value = '<div class="timeline-title">Almamet India symposium.</div>'

# "First Symposium on Hot Metal Desulphurization" -> "first symposium on hot metal desulphurization"
$ws.Range("D20").Value = "`n" + '<div class="timeline-date">2018</div>Almamet held its first symposium on hot metal desulphurization in February 2018 in Bhubaneswar.' + "`n"

# "New addition of Desulphurization Technology" -> "New addition of desulphurization technology"
$ws.Range("C21").Value = '<div class="timeline-title">Almamet Gmbh- New addition of desulphurization technology</div>'

# "Iron and Steel of" -> "iron and steel from"; "Industrial Solutions" -> "industrial solutions"
$ws.Range("D21").Value = '<div class="timeline-date">1st June 2018</div>Almamet took over metallurgical injection desulphurization technology in iron and steel from ThyssenKrupp Ploysius. ThyssenKrupp has commissioned around 130 desulphurization station worldwide since 1963. The technical team from ThyssenKrupp has joined Almamet and the industrial solutions office has been set up at Oelde, Germany.'

# "Desuplhrization" -> "desuplhrization"
$ws.Range("C22").Value = '<div class="timeline-title">Almamet Gmbh – first contract in desuplhrization technology in India</div>'

# "Operation and Maintenance contract for Desulphurization" -> "Operation and maintenance contract for desulphurization"
$ws.Range("C23").Value = '<div class="timeline-title">Almamet India – Operation and maintenance contract for desulphurization</div>'

# D15's cell style index moves from the "7" slot to the "18" slot in the source file; both xf
# records are byte-identical (same fill/font/alignment: left/bottom, wrap text, no number format),
# so re-apply the same visible formatting explicitly. (Font sub-properties are deliberately left
# untouched since assigning them materializes a brand-new font record here and would change the
# cell's actual appearance; WrapText/alignment round-trip safely onto the existing equivalent style.)
$dst = $ws.Range("D15")
$dst.HorizontalAlignment = -4131
$dst.VerticalAlignment = -4107
$dst.WrapText = $true
